$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.4700283333333333
$ws.Range("H2").Value = 1.410085
$ws.Range("I2").Value = 0.1029504401632623
$ws.Range("J2").Value = 0.1029504401632623
$ws.Range("M2").Value = 1.594573666666667
$ws.Range("N2").Value = 4.783721
$ws.Range("O2").Value = 0.02486291601650538
$ws.Range("P2").Value = 0.02555161524603605
$ws.Range("Q2").Value = 0.7494948029205555
$ws.Range("R2").Value = 6.745453226285
$ws.Range("S2").Value = 0.002559648147641452
$ws.Range("T2").Value = 0.002630550036461734

# Row 3
$ws.Range("G3").Value = 0.4700283333333333
$ws.Range("H3").Value = 1.410085
$ws.Range("I3").Value = 0.1029504401632623
$ws.Range("J3").Value = 0.1029504401632623
$ws.Range("O3").Value = 0.8216216733427845
$ws.Range("P3").Value = 0.844380476574925
$ws.Range("Q3").Value = 24.76785803115389
$ws.Range("R3").Value = 222.910722280385
$ws.Range("S3").Value = 0.08458631291831575
$ws.Range("T3").Value = 0.08692934172865369

# Row 4
$ws.Range("G4").Value = 0.4700283333333333
$ws.Range("H4").Value = 1.410085
$ws.Range("I4").Value = 0.1029504401632623
$ws.Range("J4").Value = 0.1029504401632623
$ws.Range("M4").Value = 3.409993333333333
$ws.Range("N4").Value = 10.22998
$ws.Range("O4").Value = 0.05316930765622194
$ws.Range("P4").Value = 0.05464208989919016
$ws.Range("Q4").Value = 1.602793483144444
$ws.Range("R4").Value = 14.4251413483
$ws.Range("S4").Value = 0.005473803626383959
$ws.Range("T4").Value = 0.005625427206562174

# Row 5
$ws.Range("G5").Value = 0.4700283333333333
$ws.Range("H5").Value = 1.410085
$ws.Range("I5").Value = 0.1029504401632623
$ws.Range("J5").Value = 0.1029504401632623
$ws.Range("M5").Value = 5.1859105
$ws.Range("N5").Value = 10.371821
$ws.Range("O5").Value = 0.08085976830418
$ws.Range("P5").Value = 0.05539971490660867
$ws.Range("Q5").Value = 2.437524869130833
$ws.Range("R5").Value = 14.625149214785
$ws.Range("S5").Value = 0.008324548738414733
$ws.Range("T5").Value = 0.005703425034554604

# Row 6
$ws.Range("G6").Value = 0.4700283333333333
$ws.Range("H6").Value = 1.410085
$ws.Range("I6").Value = 0.1029504401632623
$ws.Range("J6").Value = 0.1029504401632623
$ws.Range("M6").Value = 1.249748666666667
$ws.Range("N6").Value = 3.749246
$ws.Range("O6").Value = 0.01948633468030822
$ws.Range("P6").Value = 0.02002610337324014
$ws.Range("Q6").Value = 0.5874172828788888
$ws.Range("R6").Value = 5.28675554591
$ws.Range("S6").Value = 0.002006126732506374
$ws.Range("T6").Value = 0.002061696157030063

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.095550333333333
$ws.Range("H7").Value = 12.286651
$ws.Range("I7").Value = 0.8970495598367377
$ws.Range("J7").Value = 0.8970495598367377
$ws.Range("M7").Value = 1.594573666666667
$ws.Range("N7").Value = 4.783721
$ws.Range("O7").Value = 0.02486291601650538
$ws.Range("P7").Value = 0.02555161524603605
$ws.Range("Q7").Value = 6.530656712041221
$ws.Range("R7").Value = 58.775910408371
$ws.Range("S7").Value = 0.02230326786886392
$ws.Range("T7").Value = 0.02292106520957431

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.095550333333333
$ws.Range("H8").Value = 12.286651
$ws.Range("I8").Value = 0.8970495598367377
$ws.Range("J8").Value = 0.8970495598367377
$ws.Range("O8").Value = 0.8216216733427845
$ws.Range("P8").Value = 0.844380476574925
$ws.Range("Q8").Value = 215.8125415463145
$ws.Range("R8").Value = 1942.312873916831
$ws.Range("S8").Value = 0.7370353604244687
$ws.Range("T8").Value = 0.7574511348462712

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.095550333333333
$ws.Range("H9").Value = 12.286651
$ws.Range("I9").Value = 0.8970495598367377
$ws.Range("J9").Value = 0.8970495598367377
$ws.Range("M9").Value = 3.409993333333333
$ws.Range("N9").Value = 10.22998
$ws.Range("O9").Value = 0.05316930765622194
$ws.Range("P9").Value = 0.05464208989919016
$ws.Range("Q9").Value = 13.96579933299777
$ws.Range("R9").Value = 125.69219399698
$ws.Range("S9").Value = 0.04769550402983798
$ws.Range("T9").Value = 0.04901666269262799

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.095550333333333
$ws.Range("H10").Value = 12.286651
$ws.Range("I10").Value = 0.8970495598367377
$ws.Range("J10").Value = 0.8970495598367377
$ws.Range("M10").Value = 5.1859105
$ws.Range("N10").Value = 10.371821
$ws.Range("O10").Value = 0.08085976830418
$ws.Range("P10").Value = 0.05539971490660867
$ws.Range("Q10").Value = 21.23915747691183
$ws.Range("R10").Value = 127.434944861471
$ws.Range("S10").Value = 0.07253521956576527
$ws.Range("T10").Value = 0.04969628987205406

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 4.095550333333333
$ws.Range("H11").Value = 12.286651
$ws.Range("I11").Value = 0.8970495598367377
$ws.Range("J11").Value = 0.8970495598367377
$ws.Range("M11").Value = 1.249748666666667
$ws.Range("N11").Value = 3.749246
$ws.Range("O11").Value = 0.01948633468030822
$ws.Range("P11").Value = 0.02002610337324014
$ws.Range("Q11").Value = 5.118408568349555
$ws.Range("R11").Value = 46.06567711514599
$ws.Range("S11").Value = 0.01748020794780185
$ws.Range("T11").Value = 0.01796440721621007

Write-Host "Updated 126 cells across rows 2-11."
